# "Re-Run Stats + Tidy Up DiCE/SHAP"
# Adds a "Mean" summary row (row 25) below the existing results table
# (rows 3:22) with =AVERAGE(...) formulas for each metric column, and
# updates the active selection to reflect where the user left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Mean" row -------------------------------------------------
$ws.Range("B25").Value = "Mean"

# C25 gets its own (non-shared) formula; D25:G25 are entered together
# so Excel records them as one shared-formula group (t="shared").
$ws.Range("C25").Formula = "=AVERAGE(C3:C22)"
$ws.Range("D25:G25").Formula = "=AVERAGE(D3:D22)"

# Match the "0.0000" number format used by the rest of the numeric
# columns in the table.
$ws.Range("C25:G25").NumberFormat = "0.0000"

# --- Misc view state --------------------------------------------------
$ws.Range("E31").Select()

$excel.Calculate()
